# Update "想去人数" (wanted-to-go count) values in the two data sheets
# that hold the full event listing: "展览" (sheet1) and "全部类型" (sheet4).
# "演出" and "本地生活" only contain a header row, so nothing to change there.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 8317
    $ws.Range("F3").Value = 7745
    $ws.Range("F10").Value = 160
    $ws.Range("F11").Value = 228
    $ws.Range("F12").Value = 702
    $ws.Range("F14").Value = 1306
    $ws.Range("F16").Value = 50
    $ws.Range("F19").Value = 115
}
